$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.80091392993927
$ws.Range("B1").Value = 6.338734149932861
$ws.Range("C1").Value = 3.423346519470215
$ws.Range("D1").Value = 1.512985944747925
$ws.Range("E1").Value = 1.064737796783447
